# Rename column headers in existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, border) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 44948.99999999999
$newSheet.Range("B2").Value = 39
$newSheet.Range("C2").Value = -46.45343013524668
$newSheet.Range("D2").Value = 121.324817490289
$newSheet.Range("A3").Value = 44983.99999999999
$newSheet.Range("B3").Value = 48
$newSheet.Range("C3").Value = -30.40642443523056
$newSheet.Range("D3").Value = 133.7670370666816
$newSheet.Range("A4").Value = 44990.99999999999
$newSheet.Range("B4").Value = 50
$newSheet.Range("C4").Value = -26.94941568464015
$newSheet.Range("D4").Value = 133.7378242808812
$newSheet.Range("A5").Value = 44997.99999999999
$newSheet.Range("B5").Value = 52
$newSheet.Range("C5").Value = -27.27144750567446
$newSheet.Range("D5").Value = 126.3464150404987
$newSheet.Range("A6").Value = 45004.99999999999
$newSheet.Range("B6").Value = 54
$newSheet.Range("C6").Value = -24.4406487159257
$newSheet.Range("D6").Value = 127.7388844890763
$newSheet.Range("A7").Value = 45011.99999999999
$newSheet.Range("B7").Value = 56
$newSheet.Range("C7").Value = -25.2834292058949
$newSheet.Range("D7").Value = 131.8771412837678
$newSheet.Range("A8").Value = 45018.99999999999
$newSheet.Range("B8").Value = 58
$newSheet.Range("C8").Value = -25.51726211323676
$newSheet.Range("D8").Value = 138.8884227902916
$newSheet.Range("A9").Value = 45025.99999999999
$newSheet.Range("B9").Value = 60
$newSheet.Range("C9").Value = -15.78042405561499
$newSheet.Range("D9").Value = 136.992616000641
$newSheet.Range("A10").Value = 45067.99999999999
$newSheet.Range("B10").Value = 71
$newSheet.Range("C10").Value = -5.540878373315901
$newSheet.Range("D10").Value = 153.7950362003048
$newSheet.Range("A11").Value = 45074.99999999999
$newSheet.Range("B11").Value = 73
$newSheet.Range("C11").Value = -8.849736854542227
$newSheet.Range("D11").Value = 157.1299779380985
$newSheet.Range("A12").Value = 45081.99999999999
$newSheet.Range("B12").Value = 75
$newSheet.Range("C12").Value = -7.341800483022892
$newSheet.Range("D12").Value = 155.7613358281581
$newSheet.Range("A13").Value = 45088.99999999999
$newSheet.Range("B13").Value = 77
$newSheet.Range("C13").Value = -6.833572977095879
$newSheet.Range("D13").Value = 155.2339386777983
$newSheet.Range("A14").Value = 45109.99999999999
$newSheet.Range("B14").Value = 82
$newSheet.Range("C14").Value = 8.531857654634941
$newSheet.Range("D14").Value = 167.5221561203452
$newSheet.Range("A15").Value = 45130.99999999999
$newSheet.Range("B15").Value = 88
$newSheet.Range("C15").Value = 8.809188132317844
$newSheet.Range("D15").Value = 164.0182465097719
$newSheet.Range("A16").Value = 45137.99999999999
$newSheet.Range("B16").Value = 90
$newSheet.Range("C16").Value = 9.043318637181171
$newSheet.Range("D16").Value = 168.7706760900055
$newSheet.Range("A17").Value = 45144.99999999999
$newSheet.Range("B17").Value = 92
$newSheet.Range("C17").Value = 17.65175566284126
$newSheet.Range("D17").Value = 173.6495570391591
$newSheet.Range("A18").Value = 45158.99999999999
$newSheet.Range("B18").Value = 95
$newSheet.Range("C18").Value = 11.60095027780806
$newSheet.Range("D18").Value = 174.2844518145614
$newSheet.Range("A19").Value = 45165.99999999999
$newSheet.Range("B19").Value = 97
$newSheet.Range("C19").Value = 23.77352867338608
$newSheet.Range("D19").Value = 180.0049161766894
$newSheet.Range("A20").Value = 45172.99999999999
$newSheet.Range("B20").Value = 99
$newSheet.Range("C20").Value = 17.44344122696431
$newSheet.Range("D20").Value = 176.1114969602163
$newSheet.Range("A21").Value = 45179.99999999999
$newSheet.Range("B21").Value = 101
$newSheet.Range("C21").Value = 21.50304971900808
$newSheet.Range("D21").Value = 180.1195064074922
$newSheet.Range("A22").Value = 45193.99999999999
$newSheet.Range("B22").Value = 105
$newSheet.Range("C22").Value = 23.56259930817894
$newSheet.Range("D22").Value = 188.176350950537
$newSheet.Range("A23").Value = 45207.99999999999
$newSheet.Range("B23").Value = 109
$newSheet.Range("C23").Value = 26.09450524392335
$newSheet.Range("D23").Value = 183.9192244687347
$newSheet.Range("A24").Value = 45214.99999999999
$newSheet.Range("B24").Value = 110
$newSheet.Range("C24").Value = 26.22260202818235
$newSheet.Range("D24").Value = 195.0138537774238
$newSheet.Range("A25").Value = 45221.99999999999
$newSheet.Range("B25").Value = 112
$newSheet.Range("C25").Value = 37.90447625841316
$newSheet.Range("D25").Value = 192.6973703780392
$newSheet.Range("A26").Value = 45228.99999999999
$newSheet.Range("B26").Value = 114
$newSheet.Range("C26").Value = 34.0887528336957
$newSheet.Range("D26").Value = 193.5448161811777
$newSheet.Range("A27").Value = 45235.99999999999
$newSheet.Range("B27").Value = 116
$newSheet.Range("C27").Value = 35.19975569925397
$newSheet.Range("D27").Value = 192.4184354637936
$newSheet.Range("A28").Value = 45242.99999999999
$newSheet.Range("B28").Value = 118
$newSheet.Range("C28").Value = 33.21025014740063
$newSheet.Range("D28").Value = 202.0327424373133
$newSheet.Range("A29").Value = 45249.99999999999
$newSheet.Range("B29").Value = 120
$newSheet.Range("C29").Value = 44.67332221899236
$newSheet.Range("D29").Value = 201.5538942124331
$newSheet.Range("A30").Value = 45256.99999999999
$newSheet.Range("B30").Value = 122
$newSheet.Range("C30").Value = 40.8664992873058
$newSheet.Range("D30").Value = 199.2061691077515
$newSheet.Range("A31").Value = 45263.99999999999
$newSheet.Range("B31").Value = 124
$newSheet.Range("C31").Value = 48.47093844451941
$newSheet.Range("D31").Value = 203.1660776359511

# Match the date-time formatting used for the date column on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A31").PasteSpecial(-4122)
